$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-24 used to carry sequential "Processi" numbers (9-23) in column B
# with blank placeholders in column C. The refreshed keyword list no longer
# tracks those rows as processi, so drop that leftover data before rewriting
# column D with the expanded keyword list.
$ws.Range("B10:C24").ClearContents()

# Columns A and B hold small integers that must stay text (as in the rest of
# this all-text sheet), so force a text number format before assigning them -
# otherwise Excel would silently convert "7" etc. into real numbers.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2:B9").NumberFormat = "@"

# Column A
$ws.Range("A1").Value = "NumeroProcesso"
$ws.Range("A2").Value = "7"

# Column B
$ws.Range("B1").Value = "Processi"
$ws.Range("B2").Value = "1"
$ws.Range("B3").Value = "2"
$ws.Range("B4").Value = "3"
$ws.Range("B5").Value = "4"
$ws.Range("B6").Value = "5"
$ws.Range("B7").Value = "6"
$ws.Range("B8").Value = "7"
$ws.Range("B9").Value = "8"

# Column C (Ricerca) is unchanged except for the header, kept here for clarity
$ws.Range("C1").Value = "Ricerca"
$ws.Range("C2").Value = "bamboo toilet paper 5 ply 50m bamboo core 100 percent bamboo pulp plastic free FSC Ecolabel OEM"
$ws.Range("C3").Value = "bamboo jumbo tissue roll large and mini jumbo 100 percent bamboo pulp plastic free FSC OEM"
$ws.Range("C4").Value = "bamboo paper hand towels roll or folded 100 percent bamboo pulp plastic free FSC OEM"
$ws.Range("C5").Value = "A4 copy paper 80gsm 100 percent bamboo pulp plastic free FSC Ecolabel OEM custom logo"
$ws.Range("C6").Value = "notebooks and bloc-notes bamboo paper pages kraft cover plastic free FSC custom logo"
$ws.Range("C7").Value = "facial tissues 100 percent bamboo pulp pocket or box plastic free FSC Ecolabel OEM"
$ws.Range("C8").Value = "kraft paper tape water-activated gummed biodegradable plastic free FSC custom print"
$ws.Range("C9").Value = "bamboo kraft recycled paper packaging boxes and mailers plastic free FSC custom branding"

# Column D - expanded / reordered keyword list (ParoleChiave)
$ws.Range("D1").Value = "ParoleChiave"
$ws.Range("D2").Value = "packaging sostenibile"
$ws.Range("D3").Value = "imballaggio sostenibile"
$ws.Range("D4").Value = "packaging ecologico"
$ws.Range("D5").Value = "imballaggio ecologico"
$ws.Range("D6").Value = "packaging biodegradabile"
$ws.Range("D7").Value = "imballaggio biodegradabile"
$ws.Range("D8").Value = "packaging compostabile"
$ws.Range("D9").Value = "imballaggio compostabile"
$ws.Range("D10").Value = "packaging riciclabile"
$ws.Range("D11").Value = "imballaggio riciclabile"
$ws.Range("D12").Value = "carta kraft"
$ws.Range("D13").Value = "carta riciclata"
$ws.Range("D14").Value = "cellulosa di bambù"
$ws.Range("D15").Value = "fibra di bambù"
$ws.Range("D16").Value = "materiale riciclato"
$ws.Range("D17").Value = "materiale ecologico"
$ws.Range("D18").Value = "materiale sostenibile"
$ws.Range("D19").Value = "bambù naturale"
$ws.Range("D20").Value = "cartone riciclato"
$ws.Range("D21").Value = "eco friendly"
$ws.Range("D22").Value = "prodotto ecologico"
$ws.Range("D23").Value = "scatola ecologica"
$ws.Range("D24").Value = "scatola sostenibile"
$ws.Range("D25").Value = "packaging personalizzato"
$ws.Range("D26").Value = "imballaggio personalizzato"
$ws.Range("D27").Value = "stampa personalizzata"
$ws.Range("D28").Value = "etichetta ecologica"
$ws.Range("D29").Value = "busta compostabile"
$ws.Range("D30").Value = "sacchetto biodegradabile"
$ws.Range("D31").Value = "spedizione campioni"
$ws.Range("D32").Value = "bambu"
$ws.Range("D33").Value = "bambù"
$ws.Range("D34").Value = "kraft"

# Column E - new "Nome" column
$ws.Range("E1").Value = "Nome"
$ws.Range("E2").Value = "Carta"
$ws.Range("E3").Value = "Bamboo"
$ws.Range("E4").Value = "Paper"
$ws.Range("E5").Value = "A4"
$ws.Range("E6").Value = "notebooks"
$ws.Range("E7").Value = "facial tissues"
$ws.Range("E8").Value = "Kraft Paper"
$ws.Range("E9").Value = "Kraft recycled"

# Default font size bumped from 11 to 12 (workbook "Normal" style)
$wb.Styles.Item("Normal").Font.Size = 12
